# fix(Cotract): fix template contract addendum
#
# 1. Correct the misspelled company name "TRẤN THANH" -> "CHẤN THANH"
#    in the contract-date sentence.
# 2. Merge the "Địa chỉ: " run with the following address run into a
#    single run (they already share identical formatting, so Word's
#    Find/Execute naturally coalesces them when the replacement text
#    spans both runs).

$d = $word.ActiveDocument

$d.Content.Find.Execute("CÔNG TY TNHH DỊCH VỤ TRẤN THANH", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "CÔNG TY TNHH DỊCH VỤ CHẤN THANH", 2)

$d.Content.Find.Execute("Địa chỉ: A0.01 Đường Nguyễn Văn Linh, Phường Tân Phú, Quận 7, TPHCM", `
                         $true, $false, $false, $false, $false, $true, 1, $false, `
                         "Địa chỉ: A0.01 Đường Nguyễn Văn Linh, Phường Tân Phú, Quận 7, TPHCM", 2)
